# Adicionado Manager e feitas as devidas alterações
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")
$ws.Activate()

# Update the test user value referenced in B2 (usertest206 -> usertest208)
$ws.Range("B2").Value = "usertest208"

# Move the active selection to E5 (was F6)
$ws.Range("E5").Select()
